$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 40, shifting existing rows 40..135 down to 41..136
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new record
$ws.Cells.Item(40, 1).Value = 11
$ws.Cells.Item(40, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(40, 3).Value = "Bíobío"
$ws.Cells.Item(40, 4).Value = Get-Date -Year 2023 -Month 8 -Day 24 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(40, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(40, 5).Value = 8
$ws.Cells.Item(40, 6).Value = 100112012
$ws.Cells.Item(40, 7).Value = "Espinaca"
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 50
$ws.Cells.Item(40, 11).Value = 7000
$ws.Cells.Item(40, 12).Value = 7000
$ws.Cells.Item(40, 13).Value = 7000
$ws.Cells.Item(40, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(40, 15).Value = "Región Metropolitana"
$ws.Cells.Item(40, 16).Value = 700
$ws.Cells.Item(40, 17).Value = 10
$ws.Cells.Item(40, 18).Value = "Hortaliza"
